$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 1.4285715
$ws.Range("I11").Value = 1.4285715
$ws.Range("K11").Value = 1.4285715
$ws.Range("M11").Value = 138.5714285
$ws.Range("H12").Value = 216.33333
$ws.Range("I12").Value = 216.33333
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 216.33333
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -46.33332999999999
$ws.Range("N12").ClearContents()
$ws.Range("H15").Value = 519.7931
$ws.Range("I15").Value = 519.7931
$ws.Range("K15").Value = 1559.3793
$ws.Range("M15").Value = -1390.3793
$ws.Range("H31").Value = 683
$ws.Range("I31").Value = 1000
$ws.Range("K31").Value = 3000
$ws.Range("M31").Value = -2770
$ws.Range("H40").Value = 3998.7273
$ws.Range("J40").Value = 3998.7273
$ws.Range("L40").Value = 3998.7273
$ws.Range("N40").Value = -4348.7273
$ws.Range("H49").Value = 400
$ws.Range("I49").Value = 400
$ws.Range("K49").Value = 1200
$ws.Range("M49").Value = -1064
$ws.Range("H62").Value = 3812.125
$ws.Range("I62").Value = 3549.5
$ws.Range("J62").Value = 4600
$ws.Range("K62").Value = 3549.5
$ws.Range("L62").Value = 4600
$ws.Range("M62").Value = -2925.5
$ws.Range("N62").Value = -5848
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H65").Value = 3812.125
$ws.Range("I65").Value = 3549.5
$ws.Range("J65").Value = 4600
$ws.Range("K65").Value = 17747.5
$ws.Range("L65").Value = 23000
$ws.Range("M65").Value = -14627.5
$ws.Range("N65").Value = -29240
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H70").Value = 5561.75
$ws.Range("J70").Value = 2415
$ws.Range("L70").Value = 7245
$ws.Range("N70").Value = -7785
$ws.Range("H73").Value = 5561.75
$ws.Range("J73").Value = 2415
$ws.Range("L73").Value = 7245
$ws.Range("N73").Value = -9117
$ws.Range("H138").Value = 8336957.5
$ws.Range("J138").Value = 3348.9
$ws.Range("L138").Value = 10046.7
$ws.Range("N138").Value = -20326.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3787.5
$ws.Range("I132").Value = 3319.4443
$ws.Range("K132").Value = 9958.332900000001
$ws.Range("M132").Value = -7428.332900000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 2999
$ws.Range("I26").Value = 2999
$ws.Range("K26").Value = 2999
$ws.Range("M26").Value = -2707
$ws.Range("H94").Value = 1856
$ws.Range("I94").Value = 1856
$ws.Range("K94").Value = 1856
$ws.Range("M94").Value = -1405
$ws.Range("H99").Value = 1093.4445
$ws.Range("I99").Value = 1119.125
$ws.Range("J99").Value = 888
$ws.Range("K99").Value = 1119.125
$ws.Range("L99").Value = 888
$ws.Range("M99").Value = 378.875
$ws.Range("N99").Value = -3884

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 650.5
$ws.Range("J22").Value = 1000
$ws.Range("L22").Value = 1000
$ws.Range("N22").Value = -1700
$ws.Range("H43").Value = 30657
$ws.Range("J43").Value = 30657
$ws.Range("L43").Value = 30657
$ws.Range("N43").Value = -31025
$ws.Range("H58").Value = 9668.333000000001
$ws.Range("I58").Value = 7503
$ws.Range("K58").Value = 7503
$ws.Range("M58").Value = -7300
$ws.Range("H101").Value = 30657
$ws.Range("J101").Value = 30657
$ws.Range("L101").Value = 30657
$ws.Range("N101").Value = -37147
$ws.Range("H105").Value = 3999.6667
$ws.Range("I105").Value = 3999.5
$ws.Range("K105").Value = 3999.5
$ws.Range("M105").Value = -2252.5
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()
$ws.Range("H132").Value = 9051.1
$ws.Range("I132").Value = 5702.2
$ws.Range("J132").Value = 12400
$ws.Range("K132").Value = 17106.6
$ws.Range("L132").Value = 37200
$ws.Range("M132").Value = -14576.6
$ws.Range("N132").Value = -42260
$ws.Range("H134").Value = 3021.6667
$ws.Range("I134").Value = 3339
$ws.Range("K134").Value = 10017
$ws.Range("M134").Value = -7482
$ws.Range("H136").Value = 9668.333000000001
$ws.Range("I136").Value = 7503
$ws.Range("K136").Value = 22509
$ws.Range("M136").Value = -19959

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H81").Value = 5375
$ws.Range("I81").Value = 6666.6665
$ws.Range("K81").Value = 19999.9995
$ws.Range("M81").Value = -18876.9995
$ws.Range("H84").Value = 5375
$ws.Range("I84").Value = 6666.6665
$ws.Range("K84").Value = 59999.9985
$ws.Range("M84").Value = -54383.9985
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 18428572
$ws.Range("I11").Value = 19833334
$ws.Range("K11").Value = 19833334
$ws.Range("M11").Value = -19833195
$ws.Range("H134").Value = 50000
$ws.Range("J134").Value = 50000
$ws.Range("L134").Value = 150000
$ws.Range("N134").Value = -155070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3247.5
$ws.Range("I2").Value = 4000
$ws.Range("J2").Value = 990
$ws.Range("K2").Value = 4000
$ws.Range("L2").Value = 990
$ws.Range("M2").Value = -3888
$ws.Range("N2").Value = -1214
$ws.Range("H22").Value = 933.3333
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = 0
$ws.Range("N26").ClearContents()
$ws.Range("H27").Value = 933.3333
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H103").Value = 52777.5
$ws.Range("J103").Value = 52777.5
$ws.Range("L103").Value = 52777.5
$ws.Range("N103").Value = -55121.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 78273.336
$ws.Range("J46").Value = 78273.336
$ws.Range("L46").Value = 78273.336
$ws.Range("N46").Value = -78735.336
$ws.Range("H81").Value = 0
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H103").Value = 30400.666
$ws.Range("J103").Value = 30400.666
$ws.Range("L103").Value = 30400.666
$ws.Range("N103").Value = -32744.666
$ws.Range("H107").Value = 1391.3529
$ws.Range("I107").Value = 1341.091
$ws.Range("J107").Value = 1483.5
$ws.Range("K107").Value = 4023.273
$ws.Range("L107").Value = 4450.5
$ws.Range("M107").Value = -2103.273
$ws.Range("N107").Value = -8290.5
$ws.Range("H134").Value = 78273.336
$ws.Range("J134").Value = 78273.336
$ws.Range("L134").Value = 234820.008
$ws.Range("N134").Value = -239890.008
